$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('L2').Value = '*maa://24633 (55.26), *maa://30515 (69.0), *maa://34787 (71.43), ***maa://20792 (11.93), maa://39402 (84.38), ***maa://29083 (27.78)'
$ws.Range('L6').Value = 'maa://24839 (99.26)'
$ws.Range('AF7').Value = '*maa://26191 (68.0), *maa://36671 (71.74), *maa://42530 (55.56)'
$ws.Range('A8').Value = '更新日期：2024.11.16 13:18:21'
$ws.Range('P8').Value = 'maa://32931 (84.54), *maa://21916 (60.66), maa://23252 (92.42), maa://37496 (96.15), **maa://22759 (45.45)'
$ws.Range('D15').Value = '*maa://22743 (77.01), maa://22734 (83.76), *maa://30808 (63.93), ***maa://36048 (16.67)'
$ws.Range('H17').Value = 'maa://22430 (88.33), maa://39599 (83.87)'
$ws.Range('T19').Value = 'maa://24386 (98.92)'
$ws.Range('L20').Value = 'maa://41331 (83.12)'
$ws.Range('AF21').Value = 'maa://22524 (94.39), *maa://22432 (76.27)'
$ws.Range('L23').Value = 'maa://39756 (92.79), maa://39875 (93.1)'
$ws.Range('P23').Value = 'maa://30587 (91.67), *maa://29748 (75.59), ***maa://29785 (16.42), *maa://37566 (71.43)'
$ws.Range('D24').Value = 'maa://24368 (80.29)'
$ws.Range('X24').Value = 'maa://29988 (86.36), maa://23504 (92.95), **maa://22892 (39.86), *maa://25141 (77.42), maa://36663 (80.95), ***maa://22815 (23.08)'
$ws.Range('D25').Value = 'maa://29753 (95.06)'
$ws.Range('AB26').Value = 'maa://42235 (91.67)'
$ws.Range('T28').Value = 'maa://23263 (94.85), *maa://29765 (60.81)'
$ws.Range('X28').Value = 'maa://39929 (89.08), ***maa://39723 (14.29), maa://41749 (83.87)'
$ws.Range('AF29').Value = '*maa://24080 (69.33), ***maa://34960 (8.7), maa://42865 (88.0)'
$ws.Range('AB30').Value = 'maa://42979 (96.88)'
$ws.Range('L31').Value = 'maa://35926 (93.82), maa://36258 (80.72)'
$ws.Range('H32').Value = 'maa://21895 (97.06), maa://36667 (98.25), **maa://20793 (38.78), maa://22760 (100.0)'
$ws.Range('T32').Value = 'maa://41108 (87.5), maa://42859 (93.33), maa://41238 (94.23)'
$ws.Range('P33').Value = '*maa://21956 (79.41), maa://22730 (82.14)'
$ws.Range('L35').Value = 'maa://41296 (95.56)'
$ws.Range('AF38').Value = 'maa://36697 (85.44)'
$ws.Range('H53').Value = 'maa://32534 (93.26), **maa://32434 (34.78)'
$ws.Range('H55').Value = 'maa://32532 (92.18)'
$ws.Range('H59').Value = 'maa://27746 (83.5), maa://31270 (95.54)'
$ws.Range('H62').Value = 'maa://42981 (95.24)'
